# Updates the cryptos price/volume table with latest scraped values.
# Some rows also changed which coin occupies that rank (Hedera <-> InternetComputer,
# Quant <-> Flow), so B (Coin) and C (Link) are rewritten for those rows too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds plain-text numbers (dotted thousands separators like
# "23.978.81", trailing zeros like "22.20"). Without forcing text format first,
# Excel's COM layer "helpfully" reinterprets any new value that parses cleanly
# as a number (dropping a trailing zero like "22.20" -> "22.2", or silently
# retyping the cell from Text to Number even when the digits round-trip).
# Force "@" (text) on the whole data range up front so every write below
# lands as plain text, matching the original cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{ Row=2; D='24.032.88'; E='  -3.60%  ' }
    @{ Row=3; D='1.635.13'; E='  -3.30%  ' }
    @{ Row=4; D='1.003'; E='  -0.63%  ' }
    @{ Row=5; D='1.004'; E='  -0.33%  ' }
    @{ Row=6; D='306.49'; E='  -2.88%  ' }
    @{ Row=7; D='0.3904'; E='  -1.24%  ' }
    @{ Row=8; D='0.3816'; E='  -4.29%  ' }
    @{ Row=9; D='1.005'; E='  -0.52%  ' }
    @{ Row=10; E='  -6.39%  ' }
    @{ Row=11; D='1.331'; E='  -8.06%  ' }
    @{ Row=12; D='0.08354'; E='  -4.07%  ' }
    @{ Row=13; D='23.52'; E='  -7.60%  ' }
    @{ Row=14; D='7.006'; E='  -4.71%  ' }
    @{ Row=15; D='0.00001268'; E='  -4.91%  ' }
    @{ Row=16; D='7.389'; E='  -5.80%  ' }
    @{ Row=17; D='1.658.94'; E='  -2.54%  ' }
    @{ Row=18; D='94.79'; E='  +0.22%  ' }
    @{ Row=19; D='0.06872'; E='  -4.39%  ' }
    @{ Row=20; D='20.64'; E='  +1.44%  ' }
    @{ Row=21; D='6.849'; E='  -4.31%  ' }
    @{ Row=22; D='1.004'; E='  -0.26%  ' }
    @{ Row=23; D='13.45'; E='  -5.13%  ' }
    @{ Row=24; D='24.037.36'; E='  -3.61%  ' }
    @{ Row=25; E='  -2.43%  ' }
    @{ Row=26; D='2.658'; E='  -7.66%  ' }
    @{ Row=27; D='22.20'; E='  -4.27%  ' }
    @{ Row=28; D='157.30'; E='  -2.66%  ' }
    @{ Row=29; D='8.617'; E='  +7.13%  ' }
    @{ Row=30; D='139.24'; E='  -6.40%  ' }
    @{ Row=31; D='5.304'; E='  -12.71%  ' }
    @{ Row=32; D='2.396'; E='  -9.27%  ' }
    @{ Row=33; D='1.825.92'; E='  -3.43%  ' }
    @{ Row=34; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.07936'; E='  -6.17%  ' }
    @{ Row=35; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='6.762'; E='  -3.23%  ' }
    @{ Row=36; E='  -7.22%  ' }
    @{ Row=37; D='0.2657'; E='  -6.31%  ' }
    @{ Row=38; D='0.9371'; E='  -8.48%  ' }
    @{ Row=39; D='0.09143'; E='  -5.35%  ' }
    @{ Row=40; D='9.825'; E='  -8.67%  ' }
    @{ Row=41; D='1.429'; E='  -2.46%  ' }
    @{ Row=42; D='0.7455'; E='  -7.44%  ' }
    @{ Row=43; D='12.88'; E='  -7.19%  ' }
    @{ Row=44; D='15.69'; E='  -6.80%  ' }
    @{ Row=45; D='0.6804'; E='  -5.90%  ' }
    @{ Row=46; D='2.440'; E='  -6.72%  ' }
    @{ Row=47; D='4.078'; E='  -3.46%  ' }
    @{ Row=48; E='  -0.30%  ' }
    @{ Row=49; D='0.08305'; E='  -6.39%  ' }
    @{ Row=50; B='Flow'; C='https://coinranking.com/coin/QQ0NCmjVq+flow-flow'; D='1.240'; E='  -10.06%  ' }
    @{ Row=51; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='131.13'; E='  -5.06%  ' }
)

foreach ($u in $updates) {
    foreach ($col in @("B", "C", "D", "E")) {
        if ($u.ContainsKey($col)) {
            $ws.Range("$col$($u.Row)").Value = $u[$col]
        }
    }
}
